$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the outlier legend note in A3 with the expanded wording (adds AI/BI outlier
# codes and switches the separators from commas to semicolons), and bump the font
# to match the emphasis used when the new codes were introduced.
$newOutlierNote = "#Outlier: Beetle or ant crushed/Beetle or ant escaped from well/Ant not coated in any CHCs (Y); Beetle chomped (B); ant chomped (A); ant injured from the beginning(AI); beetle injured/dessicates early(BI); normal run (N)"
$a3 = $ws.Range("A3")
$a3.Value = $newOutlierNote
$a3.Font.Size = 12
$a3.Font.Color = 0
$ws.Rows.Item(3).RowHeight = 16

# Fill in the Outlier column (N) per well, now that outliers have been reviewed and
# classified so they can be plotted:
#   Wells 1,2        -> B  (beetle chomped)
#   Wells 3,4,5,6,7  -> N  (normal run)
#   Well 8           -> BI (beetle injured/dessicates early)
$outlierByWell = @{ 1 = "B"; 2 = "B"; 3 = "N"; 4 = "N"; 5 = "N"; 6 = "N"; 7 = "N"; 8 = "BI" }

# Process wells in this order (1,2 then 3-7 then 8) so that new label strings are
# created in the same order Excel would have written them as the codes were typed in.
foreach ($well in 1, 2, 3, 4, 5, 6, 7, 8) {
    $label = $outlierByWell[$well]
    for ($row = 5; $row -le 52; $row++) {
        if ([int]$ws.Cells.Item($row, 2).Value() -eq $well) {
            $ws.Cells.Item($row, 14).Value = $label
        }
    }
}

$ws.Range("N11").Select() | Out-Null
